$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the "Posted At" column keeps storing plain text dates (not Excel date serials)
$ws.Range("F2:F21").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = 'DevOps Engineer III'
$ws.Range("B2").Value = 'Cleerly'
$ws.Range("C2").Value = 'Denver, CO, US USA'
$ws.Range("D2").Value = 16.7
$ws.Range("E2").Value = 'RAG, S3, EC2, Docker, Kubernetes, CI/CD, GitHub Actions, Terraform, Git, PostgreSQL'
$ws.Range("F2").Value = '2026-02-27'
$ws.Range("G2").Value = 'https://www.indeed.com/viewjob?jk=e56f4ee4c4c94216'

# Row 3
$ws.Range("A3").Value = 'Senior Data Scientist - Product'
$ws.Range("B3").Value = 'Art of Problem Solving Academy'
$ws.Range("C3").Value = 'San Diego, CA, US USA'
$ws.Range("D3").Value = 16.7
$ws.Range("E3").Value = 'Data Scientist, Redshift, BigQuery, Git, Snowflake, BigQuery, Redshift, PySpark, Polars, Dask'
$ws.Range("F3").Value = '2026-02-27'
$ws.Range("G3").Value = 'https://www.indeed.com/viewjob?jk=60462ec1b41e6c17'

# Row 4
$ws.Range("A4").Value = 'Software Engineer II'
$ws.Range("B4").Value = 'Availity, LLC.'
$ws.Range("C4").Value = 'Remote, US USA'
$ws.Range("D4").Value = 15.6
$ws.Range("E4").Value = 'RAG, S3, EC2, Docker, Kubernetes, CI/CD, Terraform, Git, PostgreSQL, MySQL'
$ws.Range("F4").Value = '2026-02-27'
$ws.Range("G4").Value = 'https://www.indeed.com/viewjob?jk=21fc3ea476dd4728'

# Row 5
$ws.Range("A5").Value = 'Data Engineer'
$ws.Range("B5").Value = 'Health-E Commerce'
$ws.Range("C5").Value = 'Remote, US USA'
$ws.Range("D5").Value = 14.4
$ws.Range("E5").Value = 'Data Scientist, BigQuery, Synapse, Apache Airflow, Git, Snowflake, Databricks, BigQuery, Python, SQL'
$ws.Range("F5").Value = '2026-02-27'
$ws.Range("G5").Value = 'https://www.indeed.com/viewjob?jk=ec27c4fa8782043e'

# Row 6
$ws.Range("A6").Value = 'Databricks Senior Data Engineer'
$ws.Range("B6").Value = 'OZ Digital'
$ws.Range("C6").Value = 'Boca Raton, FL, US USA'
$ws.Range("D6").Value = 13.3
$ws.Range("E6").Value = 'RAG, CI/CD, Git, Snowflake, Databricks, PySpark, Kafka, Python, SQL, R'
$ws.Range("F6").Value = '2026-02-27'
$ws.Range("G6").Value = 'https://www.indeed.com/viewjob?jk=86a35b352ed2af92'

# Row 7
$ws.Range("A7").Value = 'Machine Learning Engineer'
$ws.Range("B7").Value = 'Interwell Health'
$ws.Range("C7").Value = 'Remote, US USA'
$ws.Range("D7").Value = 13.3
$ws.Range("E7").Value = 'Machine Learning Engineer, RAG, Prompt Engineering, Data Lake, AKS, CI/CD, Databricks, Python, SQL, R'
$ws.Range("F7").Value = '2026-02-27'
$ws.Range("G7").Value = 'https://www.indeed.com/viewjob?jk=888cf31c19138b41'

# Row 8
$ws.Range("A8").Value = 'Senior Quality Software Engineer'
$ws.Range("B8").Value = 'Transamerica'
$ws.Range("C8").Value = 'Denver, CO, US USA'
$ws.Range("D8").Value = 13.3
$ws.Range("E8").Value = 'RAG, S3, EC2, Docker, Kubernetes, CI/CD, Jenkins, Git, SQL, R'
$ws.Range("F8").Value = '2026-02-27'
$ws.Range("G8").Value = 'https://www.indeed.com/viewjob?jk=a424e0d64dae1bfe'

# Row 9
$ws.Range("A9").Value = 'Software Engineer – Full Stack .NET / AI Developer'
$ws.Range("B9").Value = 'nan'
$ws.Range("C9").Value = 'Bonita Springs, FL, US USA'
$ws.Range("D9").Value = 13.3
$ws.Range("E9").Value = 'RAG, Hugging Face, FAISS, Pinecone, Prompt Engineering, TensorFlow, PyTorch, CI/CD, Git, SQL'
$ws.Range("F9").Value = '2026-02-27'
$ws.Range("G9").Value = 'https://www.indeed.com/viewjob?jk=348a16c449f4600c'

# Row 10
$ws.Range("A10").Value = 'AI Software Engineer'
$ws.Range("B10").Value = 'Buyers Edge Platform'
$ws.Range("C10").Value = 'Remote, US USA'
$ws.Range("D10").Value = 12.2
$ws.Range("E10").Value = 'Data Scientist, RAG, Docker, CI/CD, Git, PostgreSQL, MySQL, Python, SQL, R'
$ws.Range("F10").Value = '2026-02-27'
$ws.Range("G10").Value = 'https://www.indeed.com/viewjob?jk=dec17c54c6cf53d9'

# Row 11
$ws.Range("A11").Value = 'Senior Software Engineer, Member AI Features'
$ws.Range("B11").Value = 'SoFi'
$ws.Range("C11").Value = 'San Francisco, CA, US USA'
$ws.Range("D11").Value = 12.2
$ws.Range("E11").Value = 'RAG, Docker, Kubernetes, CI/CD, Git, Kafka, PostgreSQL, SQL, R, Java'
$ws.Range("F11").Value = '2026-02-27'
$ws.Range("G11").Value = 'https://www.indeed.com/viewjob?jk=bfffdbe40afac374'

# Row 12
$ws.Range("A12").Value = 'Data Engineer'
$ws.Range("B12").Value = 'SoFi'
$ws.Range("C12").Value = 'Frisco, TX, US USA'
$ws.Range("D12").Value = 12.2
$ws.Range("E12").Value = 'Data Scientist, RAG, Redshift, Git, Snowflake, Redshift, PostgreSQL, Python, SQL, R'
$ws.Range("F12").Value = '2026-02-27'
$ws.Range("G12").Value = 'https://www.indeed.com/viewjob?jk=991ff35a7c2f8bc5'

# Row 13
$ws.Range("A13").Value = 'Senior Data Engineer'
$ws.Range("B13").Value = 'SoFi'
$ws.Range("C13").Value = 'Frisco, TX, US USA'
$ws.Range("D13").Value = 12.2
$ws.Range("E13").Value = 'Data Scientist, RAG, Redshift, Git, Snowflake, Redshift, PostgreSQL, Python, SQL, R'
$ws.Range("F13").Value = '2026-02-27'
$ws.Range("G13").Value = 'https://www.indeed.com/viewjob?jk=8b899c29e2338bfc'

# Row 14
$ws.Range("A14").Value = 'AI Enablement Data Engineer'
$ws.Range("B14").Value = 'IDEXX Laboratories'
$ws.Range("C14").Value = 'Westbrook, ME, US USA'
$ws.Range("D14").Value = 12.2
$ws.Range("E14").Value = 'Data Scientist, RAG, Data Lake, CI/CD, Git, Snowflake, Databricks, PySpark, Python, SQL'
$ws.Range("F14").Value = '2026-02-27'
$ws.Range("G14").Value = 'https://www.indeed.com/viewjob?jk=eed16d28c472ffad'

# Row 15
$ws.Range("A15").Value = 'Consultant Analytical Engineer Expanse'
$ws.Range("B15").Value = 'HCA Healthcare'
$ws.Range("C15").Value = 'Nashville, TN, US USA'
$ws.Range("D15").Value = 11.1
$ws.Range("E15").Value = 'RAG, Redshift, Redshift, Hadoop, Tableau, Power BI, Python, SQL, R, Scala'
$ws.Range("F15").Value = '2026-02-27'
$ws.Range("G15").Value = 'https://www.indeed.com/viewjob?jk=b0fea71c4ade6bd2'

# Row 16
$ws.Range("A16").Value = 'Sr. Analyst, Data & Research'
$ws.Range("B16").Value = 'Warner Music Group'
$ws.Range("C16").Value = 'Nashville, TN, US USA'
$ws.Range("D16").Value = 11.1
$ws.Range("E16").Value = 'RAG, BigQuery, Git, Snowflake, BigQuery, Tableau, Python, SQL, R, Scala'
$ws.Range("F16").Value = '2026-02-27'
$ws.Range("G16").Value = 'https://www.indeed.com/viewjob?jk=2c05d60c76967e19'

# Row 17
$ws.Range("A17").Value = 'AI Engineer - Agent Automation'
$ws.Range("B17").Value = 'Zoom Communications'
$ws.Range("C17").Value = 'San Jose, CA, US USA'
$ws.Range("D17").Value = 11.1
$ws.Range("E17").Value = 'AI Engineer, LangChain, RAG, Prompt Engineering, PyTorch, Docker, Kubernetes, Python, R, Java'
$ws.Range("F17").Value = '2026-02-27'
$ws.Range("G17").Value = 'https://www.indeed.com/viewjob?jk=f5af786d5985c54f'

# Row 18
$ws.Range("A18").Value = 'Sr. AI Software Engineer'
$ws.Range("B18").Value = 'Technoviz LLC'
$ws.Range("C18").Value = 'US USA'
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = 'LangChain, RAG, Hugging Face, Docker, Kubernetes, Git, R, Scala, Optimization'
$ws.Range("F18").Value = '2026-02-27'
$ws.Range("G18").Value = 'https://www.indeed.com/viewjob?jk=c72678b7c20061ce'

# Row 19
$ws.Range("A19").Value = 'GTM Analytics Engineer'
$ws.Range("B19").Value = 'Gusto'
$ws.Range("C19").Value = 'New York, NY, US USA'
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 'Redshift, Git, Snowflake, Redshift, Tableau, Python, SQL, R, Optimization'
$ws.Range("F19").Value = '2026-02-27'
$ws.Range("G19").Value = 'https://www.indeed.com/viewjob?jk=ba485f3b26852488'

# Row 20
$ws.Range("A20").Value = 'Data Scientist'
$ws.Range("B20").Value = 'The Home Depot'
$ws.Range("C20").Value = 'Atlanta, GA, US USA'
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = 'Data Scientist, RAG, BigQuery, BigQuery, Tableau, Python, SQL, R, Optimization'
$ws.Range("F20").Value = '2026-02-27'
$ws.Range("G20").Value = 'https://www.indeed.com/viewjob?jk=f0dffbe9c159b120'

# Row 21
$ws.Range("A21").Value = 'AI / ML Engineer'
$ws.Range("B21").Value = 'Third Way Health, Inc.'
$ws.Range("C21").Value = 'Cambridge, MA, US USA'
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = 'Data Scientist, RAG, PyTorch, Kinesis, MLflow, CI/CD, Kafka, Python, R'
$ws.Range("F21").Value = '2026-02-27'
$ws.Range("G21").Value = 'https://www.indeed.com/viewjob?jk=01b8d0cbb789aae6'
